$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks so we can rebuild them cleanly for the new layout
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range('A2').Value = '2026-02-14 01:51:54'
$ws.Range('B2').Value = '【基礎エンジニア歓迎】暗号資産取引所APIを使ったPython自動化ツール開発'
$ws.Range('C2').Value = 'システム開発'
$ws.Range('D2').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E2').Value = '期限情報なし'
$ws.Range('F2').Value = 'https://www.lancers.jp/work/detail/5491124'
$ws.Range('G2').Value = 543
$ws.Range('H2').Value = '🔥Python,API ◆ツール,開発'

# Row 3
$ws.Range('A3').Value = '2026-02-14 01:51:54'
$ws.Range('B3').Value = '最新AI活用、書き伝票から在庫更新請求入金消込までの完全自動化スキーム構築Claude/Gemini'
$ws.Range('C3').Value = 'システム開発'
$ws.Range('D3').Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range('E3').Value = '期限情報なし'
$ws.Range('F3').Value = 'https://www.lancers.jp/work/detail/5490911'
$ws.Range('G3').Value = 395
$ws.Range('H3').Value = '🔥AI,Ai ◆自動化'

# Row 4
$ws.Range('A4').Value = '2026-02-14 01:51:54'
$ws.Range('B4').Value = '産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)'
$ws.Range('C4').Value = 'システム開発'
$ws.Range('D4').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E4').Value = '期限情報なし'
$ws.Range('F4').Value = 'https://www.lancers.jp/work/detail/5450864'
$ws.Range('G4').Value = 383
$ws.Range('H4').Value = '🔥AI,Ai ◆開発'

# Row 5
$ws.Range('A5').Value = '2026-02-14 01:51:54'
$ws.Range('B5').Value = '【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集'
$ws.Range('C5').Value = 'システム開発'
$ws.Range('D5').Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range('E5').Value = '期限情報なし'
$ws.Range('F5').Value = 'https://www.lancers.jp/work/detail/5460294'
$ws.Range('G5').Value = 375
$ws.Range('H5').Value = '🔥AI,Ai ◆開発'

# Row 6
$ws.Range('A6').Value = '2026-02-14 01:51:54'
$ws.Range('B6').Value = '【急募】ビジネス向けAIエージェント開発支援のパートナー募集'
$ws.Range('C6').Value = 'システム開発'
$ws.Range('D6').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E6').Value = '期限情報なし'
$ws.Range('F6').Value = 'https://www.lancers.jp/work/detail/5490828'
$ws.Range('G6').Value = 368
$ws.Range('H6').Value = '🔥AI,Ai ◆開発'

# Row 7
$ws.Range('A7').Value = '2026-02-14 01:51:54'
$ws.Range('B7').Value = '【急募・フルリモート】React Native アプリ開発エンジニア募集'
$ws.Range('C7').Value = 'システム開発'
$ws.Range('D7').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E7').Value = '期限情報なし'
$ws.Range('F7').Value = 'https://www.lancers.jp/work/detail/5491190'
$ws.Range('G7').Value = 218
$ws.Range('H7').Value = '🔥React ◆開発 ◇アプリ'

# Row 8
$ws.Range('A8').Value = '2026-02-14 01:51:54'
$ws.Range('B8').Value = '【医療機関向け業務改善サービスの新規開発】WEBアプリ開発におけるフルスタック開発担当者募集'
$ws.Range('C8').Value = 'システム開発'
$ws.Range('D8').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E8').Value = '期限情報なし'
$ws.Range('F8').Value = 'https://www.lancers.jp/work/detail/5473940'
$ws.Range('G8').Value = 135
$ws.Range('H8').Value = '◆開発 ◇業務改善'

# Row 9
$ws.Range('A9').Value = '2026-02-14 01:51:54'
$ws.Range('B9').Value = '【急募】アンドロイドタブレット向け将棋アプリ開発者募集'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('F9').Value = 'https://www.lancers.jp/work/detail/5491505'
$ws.Range('G9').Value = 88
$ws.Range('H9').Value = '◆開発 ◇アプリ'

# Row 10
$ws.Range('A10').Value = '2026-02-14 01:51:54'
$ws.Range('B10').Value = 'PHP/Laravelエンジニア募集(大規模Webシステム/フルリモート)'
$ws.Range('C10').Value = 'システム開発'
$ws.Range('D10').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E10').Value = '期限情報なし'
$ws.Range('F10').Value = 'https://www.lancers.jp/work/detail/5490679'
$ws.Range('G10').Value = 55
$ws.Range('H10').Value = '○PHP'

# Row 11
$ws.Range('A11').Value = '2026-02-14 01:51:54'
$ws.Range('B11').Value = 'bubbleで構築したサイトの修正対応'
$ws.Range('C11').Value = 'システム開発'
$ws.Range('D11').Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range('E11').Value = '期限情報なし'
$ws.Range('F11').Value = 'https://www.lancers.jp/work/detail/5491578'
$ws.Range('G11').Value = 30
$ws.Range('H11').Value = '◇サイト'

# Row 12
$ws.Range('A12').Value = '2026-02-14 01:51:54'
$ws.Range('B12').Value = 'bubbleで構築したサイトの修正対応'
$ws.Range('C12').Value = 'システム開発'
$ws.Range('D12').Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range('E12').Value = '期限情報なし'
$ws.Range('F12').Value = 'https://www.lancers.jp/work/detail/5491569'
$ws.Range('G12').Value = 30
$ws.Range('H12').Value = '◇サイト'

# Row 13
$ws.Range('A13').Value = '2026-02-14 01:51:54'
$ws.Range('B13').Value = '【急募】1週間でGASを用いたカレンダー同期システム構築'
$ws.Range('C13').Value = 'システム開発'
$ws.Range('D13').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E13').Value = '期限情報なし'
$ws.Range('F13').Value = 'https://www.lancers.jp/work/detail/5491203'
$ws.Range('G13').Value = 33
$ws.Range('H13').Value = ""

# Row 14
$ws.Range('A14').Value = '2026-02-14 01:51:54'
$ws.Range('B14').Value = 'yahooプレイス用Worepressプラグイン 投稿記事をyahooプレイスのお知らせに要約投稿'
$ws.Range('C14').Value = 'システム開発'
$ws.Range('D14').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E14').Value = '期限情報なし'
$ws.Range('F14').Value = 'https://www.lancers.jp/work/detail/5477871'
$ws.Range('G14').Value = 18
$ws.Range('H14').Value = ""

# Row 15
$ws.Range('A15').Value = '2026-02-14 01:51:54'
$ws.Range('B15').Value = '【RHEL5 → 新OS中継サーバ】メール基盤構築・疎通確認まで対応できる方募集'
$ws.Range('C15').Value = 'システム開発'
$ws.Range('D15').Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range('E15').Value = '期限情報なし'
$ws.Range('F15').Value = 'https://www.lancers.jp/work/detail/5491086'
$ws.Range('G15').Value = 18
$ws.Range('H15').Value = ""

# Row 16
$ws.Range('A16').Value = '2026-02-14 01:51:54'
$ws.Range('B16').Value = '【3,000円 / 急募】GitHubとVercelの連携設定エラーの解消'
$ws.Range('C16').Value = 'システム開発'
$ws.Range('D16').Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range('E16').Value = '期限情報なし'
$ws.Range('F16').Value = 'https://www.lancers.jp/work/detail/5491643'
$ws.Range('G16').Value = 10
$ws.Range('H16').Value = ""

# Row 17
$ws.Range('A17').Value = '2026-02-14 01:51:54'
$ws.Range('B17').Value = '放置中の法人ドメインを持っている企業様'
$ws.Range('C17').Value = 'システム開発'
$ws.Range('D17').Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range('E17').Value = '期限情報なし'
$ws.Range('F17').Value = 'https://www.lancers.jp/work/detail/5490905'
$ws.Range('G17').Value = 10
$ws.Range('H17').Value = ""

# Re-add hyperlinks for the URL column (F) on every data row
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5491124')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5490911')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5450864')
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5460294')
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5490828')
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5491190')
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5473940')
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5491505')
$ws.Hyperlinks.Add($ws.Range('F10'), 'https://www.lancers.jp/work/detail/5490679')
$ws.Hyperlinks.Add($ws.Range('F11'), 'https://www.lancers.jp/work/detail/5491578')
$ws.Hyperlinks.Add($ws.Range('F12'), 'https://www.lancers.jp/work/detail/5491569')
$ws.Hyperlinks.Add($ws.Range('F13'), 'https://www.lancers.jp/work/detail/5491203')
$ws.Hyperlinks.Add($ws.Range('F14'), 'https://www.lancers.jp/work/detail/5477871')
$ws.Hyperlinks.Add($ws.Range('F15'), 'https://www.lancers.jp/work/detail/5491086')
$ws.Hyperlinks.Add($ws.Range('F16'), 'https://www.lancers.jp/work/detail/5491643')
$ws.Hyperlinks.Add($ws.Range('F17'), 'https://www.lancers.jp/work/detail/5490905')
